$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (Modelos Informacionais) - prova grade entered: Unidade 4 = 60, AP/RP = "AP"
$ws.Range("J4").Value = 60
$ws.Range("L4").Value = "AP"

# Row 8 (BIM 4D - Planejamento e Controle de Obras) - prova grade entered: Unidade 4 = 54, AP/RP = "AP"
$ws.Range("J8").Value = 54
$ws.Range("L8").Value = "AP"

# K10/K12/K14 (still-empty rows) pick up the "even row" number style
# (grey font + shaded fill + 0.00 number format), matching rows 4/6/8's look.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("K10").PasteSpecial(-4122) | Out-Null
$ws.Range("K10").NumberFormat = "0.00"

$ws.Range("A10").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null
$ws.Range("K12").NumberFormat = "0.00"

$ws.Range("A10").Copy() | Out-Null
$ws.Range("K14").PasteSpecial(-4122) | Out-Null
$ws.Range("K14").NumberFormat = "0.00"

# K11/K13 pick up the "odd row" number style (grey font, no fill, 0.00 number format).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("K11").PasteSpecial(-4122) | Out-Null
$ws.Range("K11").NumberFormat = "0.00"

$ws.Range("A3").Copy() | Out-Null
$ws.Range("K13").PasteSpecial(-4122) | Out-Null
$ws.Range("K13").NumberFormat = "0.00"

# Leave the same selection behind as the saved workbook (E21:K21).
$ws.Range("E21:K21").Select() | Out-Null
